$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "2026-01-30 05:34"
$ws.Range("B9").Value = 36
$ws.Range("C9").Value = 6
